# Milestone 1 - final presentation edits
#  1. Remove the "Misc." slide from the deck (it was slide 8).
#  2. Refresh the cached footer/date-placeholder text from "2019-10-01" to
#     "2019-10-02" on the slide master and on every slide layout.

$p = $ppt.ActivePresentation

# --- 1. Delete the "Misc." slide --------------------------------------------------
$miscIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Name -like "Title*" -and $sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Misc.") {
            $miscIndex = $i
        }
    }
}
if ($miscIndex -eq -1) {
    $miscIndex = 8
}
$p.Slides.Item($miscIndex).Delete()

# --- 2. Update the cached date-placeholder text -----------------------------------
$newDate = "2019-10-02"

$m = $p.SlideMaster
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
